# Port-level landings table (Table18) update
# - Row 40 used to be a stray "Santa Cruz" port-totals row that actually
#   duplicated the "MONTEREY AREA TOTALS" label from the last area-total
#   row; it is turned into a generic "Totals" row labelled with the
#   "MONTEREY AREA TOTALS" text moved into column A.
# - Column A is widened to match column B (both columns now share the
#   same best-fit width).
# - The active selection is left on the whole of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A40: "Santa Cruz" -> "MONTEREY AREA TOTALS"
$ws.Range("A40").Value = "MONTEREY AREA TOTALS"

# B40: "MONTEREY AREA TOTALS" -> "Totals" (new shared string)
$ws.Range("B40").Value = "Totals"

# Column A now matches column B's (best-fit) width
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Select the entirety of column A, as the author left it
$ws.Columns.Item(1).Select()
